# Brooke Cooper Reflections.docx - apply commit "Uploading one change.  A turn of phrase."
#
# The diff collapses three spell-check-wrapped words (Keneau / Dartanion / C'est)
# back into their surrounding runs (removing the now-stale <w:proofErr/> markers),
# and inserts a new aside " (but not literally)" after "In a word".

$d = $word.ActiveDocument

# --- Change 1: "...picture the [Keneau] Reeves..." -------------------------
# Touch the whole sentence that spans the proofErr-wrapped "Keneau" so the
# surrounding runs (and the now-unneeded spell-check markers) collapse into one.
$r1 = $d.Content.Find.Execute(
    "We are, in effect, uploading programs into our cerebral cortex.  I can just picture the Keneau Reeves character Neo saying, " + [char]34 + "I know Python." + [char]34 + "  Some of it feels like a blur because of the course pace.  Some of it feels incomplete (and obviously is), but it's the feeling I'm focused on.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We are, in effect, uploading programs into our cerebral cortex.  I can just picture the Keneau Reeves character Neo saying, " + [char]34 + "I know Python." + [char]34 + "  Some of it feels like a blur because of the course pace.  Some of it feels incomplete (and obviously is), but it's the feeling I'm focused on.",
    2)
Write-Host "Change 1 (Keneau):" $r1

# --- Change 2: "...presenting as [Dartanion] Williams." ---------------------
$r2 = $d.Content.Find.Execute(
    "in presenting as Dartanion Williams",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "in presenting as Dartanion Williams",
    2)
Write-Host "Change 2 (Dartanion):" $r2

# --- Change 3: "In a word, more machine learning..." -> insert an aside ----
$r3 = $d.Content.Find.Execute(
    "In a word,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In a word (but not literally),",
    2)
Write-Host "Change 3 (turn of phrase):" $r3

# --- Change 4: "...after all this effort.  [C'est] la vie." ----------------
$r4 = $d.Content.Find.Execute(
    "At this point, there is a real chance I will end up doing nothing worthwhile or at least nothing interesting to me after all this effort.  C'est la vie.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "At this point, there is a real chance I will end up doing nothing worthwhile or at least nothing interesting to me after all this effort.  C'est la vie.",
    2)
Write-Host "Change 4 (C'est la vie):" $r4
